$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ----------------------------------------------------------------------
# Hunk 1: split the "RESTful web App..." paragraph (index 37) into three
# paragraphs: the rewritten intro sentence, an (empty) paragraph that now
# just carries the relocated _GoBack bookmark, and a new "Functionality
# includes..." paragraph.
# ----------------------------------------------------------------------
$target = $d.Paragraphs(37)
$target.Range.InsertParagraphBefore() | Out-Null
$target.Range.InsertParagraphBefore() | Out-Null
# After the two inserts: 37 and 38 are fresh blank paragraphs (inherited
# the old bold pPr from paragraph 37), 39 is the untouched original
# paragraph.

$para1Xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr><w:t>RESTful web App built with Python</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr><w:t>,</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr><w:t xml:space='preserve'> Flask, SQLAlchemy, Jinja2 and Bootstrap.</w:t></w:r>" +
  "</w:p>"

$para2Xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:bCs/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr></w:pPr>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
  "</w:p>"

$para3Xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:bCs/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:bCs/><w:color w:val='000000'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:eastAsia='da-DK'/></w:rPr><w:t>Functionality includes OAuth v2.0 integration for Google accounts, and CRUD with CSRF protection for entries and local permission systems.</w:t></w:r>" +
  "</w:p>"

$d.Paragraphs(37).Range.InsertXML($para1Xml) | Out-Null
$d.Paragraphs(38).Range.InsertXML($para2Xml) | Out-Null
$d.Paragraphs(39).Range.InsertXML($para3Xml) | Out-Null

# ----------------------------------------------------------------------
# Hunk 2: the _GoBack bookmark used to live inside the "Delete a
# category" route paragraph, just before the final "delete" run; it was
# relocated above, so remove it from its old location here. Two new
# paragraphs were inserted above, so this paragraph shifted from 71 -> 73.
# ----------------------------------------------------------------------
$para71Xml = "<w:p $wns>" +
  "<w:pPr><w:pStyle w:val='Listeafsnit'/><w:widowControl w:val='0'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='10'/></w:numPr><w:autoSpaceDE w:val='0'/><w:autoSpaceDN w:val='0'/><w:adjustRightInd w:val='0'/><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>/</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:i/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>category-name</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>/</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>category_id/</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>delete</w:t></w:r>" +
  "</w:p>"

$d.Paragraphs(73).Range.InsertXML($para71Xml) | Out-Null

# ----------------------------------------------------------------------
# Hunks 3-6: update the slug-style route examples so they spell out
# Flask's <converter:name> variable-rule syntax (category_name,
# category_id, item_name, item_id) instead of the old dashed placeholders.
# The affected paragraphs shifted from 75/77/79/81 -> 77/79/81/83.
# ----------------------------------------------------------------------
$listPPr = "<w:pPr><w:pStyle w:val='Listeafsnit'/><w:widowControl w:val='0'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='12'/></w:numPr><w:autoSpaceDE w:val='0'/><w:autoSpaceDN w:val='0'/><w:adjustRightInd w:val='0'/><w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>"
$plainRPr = "<w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>"
$italicRPr = "<w:rPr><w:rFonts w:ascii='Helvetica' w:hAnsi='Helvetica' w:cs='Helvetica'/><w:i/><w:color w:val='353535'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>"

$para75Xml = "<w:p $wns>" + $listPPr +
  "<w:r>$plainRPr<w:t>/</w:t></w:r>" +
  "<w:r>$italicRPr<w:t>category</w:t></w:r>" +
  "<w:r>$italicRPr<w:t>_</w:t></w:r>" +
  "<w:r>$italicRPr<w:t>name</w:t></w:r>" +
  "<w:r>$italicRPr<w:t>/category_id</w:t></w:r>" +
  "</w:p>"

$para77Xml = "<w:p $wns>" + $listPPr +
  "<w:r>$plainRPr<w:t>/category_</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>name</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>/category_id/item_</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>name/</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>item_id</w:t></w:r>" +
  "</w:p>"

$para79Xml = "<w:p $wns>" + $listPPr +
  "<w:r>$plainRPr<w:t>/</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>category_</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>name/</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>category_id/</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>new</w:t></w:r>" +
  "</w:p>"

$para81Xml = "<w:p $wns>" + $listPPr +
  "<w:r>$plainRPr<w:t>/category_</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>name/</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>category_id/item_name/item_id</w:t></w:r>" +
  "<w:r>$plainRPr<w:t>/edit</w:t></w:r>" +
  "</w:p>"

$d.Paragraphs(77).Range.InsertXML($para75Xml) | Out-Null
$d.Paragraphs(79).Range.InsertXML($para77Xml) | Out-Null
$d.Paragraphs(81).Range.InsertXML($para79Xml) | Out-Null
$d.Paragraphs(83).Range.InsertXML($para81Xml) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
